$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to remain text even when the new value parses as a
    # plain number (Excel would otherwise silently convert it to a Number).
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.180.49"
$ws.Range("E2").Value = "  +1.35%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.646.39"
$ws.Range("E3").Value = "  +0.19%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.15%  "

# Row 5 - BNB
Set-TextValue "D5" "217.02"
$ws.Range("E5").Value = "  +0.04%  "

# Row 6 - XRP
Set-TextValue "D6" "0.516"
$ws.Range("E6").Value = "  +2.41%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.12%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  +1.41%  "

# Row 9 - Dogecoin
Set-TextValue "D9" "0.0628"
$ws.Range("E9").Value = "  +1.31%  "

# Row 10 - Solana
Set-TextValue "D10" "19.95"
$ws.Range("E10").Value = "  +1.34%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.54%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.876.87"
$ws.Range("E12").Value = "  +0.23%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.640.93"
$ws.Range("E13").Value = "  -0.22%  "

# Row 14 - Polkadot
Set-TextValue "D14" "4.15"
$ws.Range("E14").Value = "  +0.45%  "

# Row 15 - Polygon
Set-TextValue "D15" "0.541"
$ws.Range("E15").Value = "  +2.70%  "

# Row 16 - Litecoin
Set-TextValue "D16" "67.62"
$ws.Range("E16").Value = "  +2.27%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "27.165.99"
$ws.Range("E17").Value = "  +1.22%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "0.0₃0739"
$ws.Range("E18").Value = "  +1.32%  "

# Row 19 - BitcoinCash
Set-TextValue "D19" "219.07"
$ws.Range("E19").Value = "  +0.62%  "

# Row 21 - Toncoin
$ws.Range("E21").Value = "  +5.27%  "

# Row 22 - Chainlink
$ws.Range("E22").Value = "  +3.03%  "

# Row 23 - Uniswap
Set-TextValue "D23" "4.41"
$ws.Range("E23").Value = "  +0.37%  "

# Row 24 - Avalanche
Set-TextValue "D24" "9.21"
$ws.Range("E24").Value = "  +0.52%  "

# Row 25 - Monero
Set-TextValue "D25" "147.92"

# Row 26 - Cosmos
Set-TextValue "D26" "7.57"
$ws.Range("E26").Value = "  +3.13%  "

# Row 27 - BinanceUSD
$ws.Range("E27").Value = "  -0.15%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  +0.08%  "

# Row 29 - EthereumClassic
Set-TextValue "D29" "15.75"
$ws.Range("E29").Value = "  -0.42%  "

# Row 30 - Hedera
Set-TextValue "D30" "0.0508"
$ws.Range("E30").Value = "  -0.69%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -0.03%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +0.70%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  +1.31%  "

# Row 34 - LidoDAOToken
Set-TextValue "D34" "1.58"
$ws.Range("E34").Value = "  +1.81%  "

# Row 35 - Maker
$ws.Range("D35").Value = "1.264.67"
$ws.Range("E35").Value = "  +1.40%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  +0.39%  "

# Row 37 - VeChain
$ws.Range("E37").Value = "  +1.65%  "

# Row 38 - ImmutableX
Set-TextValue "D38" "0.547"
$ws.Range("E38").Value = "  +2.53%  "

# Row 39 - ARBITRUM
$ws.Range("E39").Value = "  +1.84%  "

# Row 40 - PaxDollar
$ws.Range("E40").Value = "  -0.16%  "

# Row 41 - TrustWalletToken
$ws.Range("E41").Value = "  +0.18%  "

# Row 42 - MXToken
$ws.Range("E42").Value = "  +6.15%  "

# Row 43 - FraxShare
$ws.Range("E43").Value = "  +1.52%  "

# Row 44 - RocketPoolETH
$ws.Range("D44").Value = "1.786.92"
$ws.Range("E44").Value = "  +0.09%  "

# Row 45 - Aave
Set-TextValue "D45" "61.95"
$ws.Range("E45").Value = "  +1.87%  "

# Row 46 - Quant
Set-TextValue "D46" "91.45"

# Row 47 - RenderToken
$ws.Range("E47").Value = "  +0.89%  "

# Row 48 - BabyDogeCoin
$ws.Range("D48").Value = "0.0₆0107"
$ws.Range("E48").Value = "  +2.29%  "

# Row 49 - Cronos
Set-TextValue "D49" "0.0513"
$ws.Range("E49").Value = "  -0.07%  "

# Row 50 - EnergySwap
Set-TextValue "D50" "7.62"
$ws.Range("E50").Value = "  +0.95%  "

# Row 51 - Algorand
Set-TextValue "D51" "0.0973"
$ws.Range("E51").Value = "  +0.31%  "
